# "Koordinaten des Zentrums angepasst * 20 neue Taxistaende eingepflegt"
#
# TaxiStands (sheet2): rows 2-3 (the "center" stands) get corrected
# lat/long/capacity values, and 20 brand-new taxi-stand rows (4-26) are
# appended below them. The sheet becomes the active tab with a specific
# cell selected, and the data cells get an explicit (non-default) cell
# style split between the coordinate columns (A,B) and the capacity
# columns (C,D,E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaxiStands")

# --- full TaxiStands data for rows 2..26 -----------------------------
# columns: row, locationLat, locationLong, parkCapacity, chargeCapacity, outputkWh
$data = @(
    @(2, 52.386400000000002, 9.7109100000000002, 1, 100, 100),
    @(3, 52.389690000000002, 9.7207399999999993, 1, 100, 100),
    @(4, 52.381039999999999, 9.7362800000000007, 1, 100, 100),
    @(5, 52.37323, 9.72865, 1, 100, 100),
    @(6, 52.373089999999998, 9.7082700000000006, 1, 100, 100),
    @(7, 52.380740000000003, 9.6991599999999991, 1, 100, 100),
    @(8, 52.366759999999999, 9.7177699999999998, 1, 100, 100),
    @(9, 52.373089999999998, 9.7440700000000007, 1, 100, 100),
    @(10, 52.390889999999999, 9.7360600000000002, 1, 100, 100),
    @(11, 52.397959999999998, 9.7083300000000001, 2, 200, 200),
    @(12, 52.374769999999998, 9.6938300000000002, 2, 200, 200),
    @(13, 52.365789999999997, 9.7010799999999993, 2, 200, 200),
    @(14, 52.364780000000003, 9.7360000000000007, 2, 200, 200),
    @(15, 52.37377, 9.7597199999999997, 2, 200, 200),
    @(16, 52.386130000000001, 9.7498400000000007, 2, 200, 200),
    @(17, 52.399000000000001, 9.7409400000000002, 2, 200, 200),
    @(18, 52.366669999999999, 9.6840899999999994, 2, 200, 200),
    @(19, 52.359400000000001, 9.7108000000000008, 2, 200, 200),
    @(20, 52.364429999999999, 9.7525899999999996, 2, 200, 200),
    @(21, 52.381419999999999, 9.7661999999999995, 2, 200, 200),
    @(22, 52.356819999999999, 9.76, 2, 200, 200),
    @(23, 52.365400000000001, 9.7692200000000007, 2, 200, 200),
    @(24, 52.349710000000002, 9.7695500000000006, 2, 200, 200),
    @(25, 52.399050000000003, 9.7587399999999995, 2, 200, 200),
    @(26, 52.397950000000002, 9.7270000000000003, 2, 200, 200)
)

$lastRow = 26

# --- give the two column groups their own explicit cell style --------
# A:B (coordinates) vs C:E (capacities) end up as two distinct, explicit
# styles rather than the workbook default.
$ws.Range("C2:E$lastRow").NumberFormat = "0"
$ws.Range("A2:B$lastRow").NumberFormat = "General"

# --- write the values --------------------------------------------------
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# --- make TaxiStands the active sheet with the new selection ---------
$ws.Activate()
$ws.Range("I24").Select()
